$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
# G2 first so the "Succès" shared string becomes index 23, and the new
# green font/style becomes the first of the three newly added styles.
$ws.Range("G2").Value = "Succès"
$ws.Range("G2").Font.ColorIndex = 17
$ws.Range("F2").Value = "This is a demo for NORAUI (Non-Regression Automation for User Interfaces)."

# --- Row 3 ---
$ws.Range("G3").Value = "Échec : The city is Paris!!"
$ws.Range("G3").Font.ColorIndex = 10

# --- Row 4 ---
$ws.Range("G4").Value = "Échec : Saisie « Input Select field » dans demo."
$ws.Range("G4").Font.ColorIndex = 10

# --- Row 5 ---
# F5 previously carried an (empty) style; clear it back to the default
# "Normal" style before writing the value so no explicit s="" survives.
$ws.Range("F5").Style = "Normal"
$ws.Range("F5").Value = "This is a demo for NORAUI (Non-Regression Automation for User Interfaces)."
$ws.Range("G5").Value = "Succès"
$ws.Range("G5").Font.ColorIndex = 17

# --- Row 6 ---
$ws.Range("G6").Value = "Échec : Accès à l'action « no exist element » dans demo."
$ws.Range("G6").Font.ColorIndex = 10

# --- Row 7 ---
$ws.Range("G7").Value = "Échec : La donnée « city » fournie ne peut pas être vide."
$ws.Range("G7").Font.ColorIndex = 10

# --- Row 8 ---
$ws.Range("G8").Value = "Échec : La donnée « element » fournie ne peut pas être vide, la donnée « element2 » fournie ne peut pas être vide."
$ws.Range("G8").Font.ColorIndex = 10

# --- Row 9 ---
$ws.Range("F9").Value = "This is a demo for NORAUI (Non-Regression Automation for User Interfaces)."
$ws.Range("G9").Value = "Échec : Absence « -input_text_field » dans demo."
$ws.Range("G9").Font.ColorIndex = 10

# A third font/style (brown, indexed 53) was also registered in the
# workbook's style table without being applied to any visible cell;
# reproduce that by touching a scratch cell and clearing it again.
$ws.Range("ZZ1").Font.ColorIndex = 53
$ws.Range("ZZ1").Clear()
